$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 and J1, copying the style from H1 (bold, bordered, centered header style)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(6,5,6,5,6,7,6,7,7,7,5,6,8,8,9,7,9,9,8,9,8,9,7,7,7,10,6,8,5,6,6,7,9,9,10,9,6,8,9,9,8,9,9,8,10,9,9,9,9,9,8,9,9,9,9,9,9,8,9,9,9,9,8,10,9,8,9,8,5,5,4,2)
$jValues = @(6,6,7,6,7,7,6,7,7,7,5,6,8,8,9,7,9,9,8,9,8,9,7,7,7,10,6,8,5,6,6,7,9,9,10,9,6,8,9,9,8,9,9,8,10,9,9,9,9,9,8,9,9,9,9,9,9,8,9,9,9,9,9,10,9,8,9,8,5,5,4,2)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}

$excel.CutCopyMode = 0
